$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Title (row 5) changes from "Session type" to "NG-Imm Session Type VS"
$ws.Range("B5").Value = "NG-Imm Session Type VS"

# Date (row 8) changes to the new generation timestamp
$ws.Range("B8").Value = "2025-06-24T09:13:37+01:00"
